$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: remove the trailing " (old version)" run that follows the heading
# "Reconstruction algorithm #1".
# ---------------------------------------------------------------------------
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("Reconstruction algorithm #1 (old version)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $delStart1 = $rng1.Start + ("Reconstruction algorithm #1").Length
    $delRng1 = $d.Range($delStart1, $rng1.End)
    $delRng1.Delete()
}

# ---------------------------------------------------------------------------
# Edit 2: in the caption ". (B) A schematic of reconstruction algorithm #1
# (old version) based on the phase encoding scheme ..." remove " (old
# version)" and drop a (now-orphaned) "_GoBack" bookmark right between the
# "#" and the "1" (mirrors Word leaving its last-edit-position bookmark
# behind once the "(old version)" phrase was trimmed out of the sentence).
# ---------------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("schematic of reconstruction algorithm #1 (old version)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPoint = $rng2.Start + ("schematic of reconstruction algorithm #").Length
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $rng2b = $d.Content.Duplicate
    $found2b = $rng2b.Find.Execute(" (old version) based", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2b) {
        $delEnd2 = $rng2b.Start + (" (old version)").Length
        $delRng2 = $d.Range($rng2b.Start, $delEnd2)
        $delRng2.Delete()
    }
}

# ---------------------------------------------------------------------------
# Edit 3: the Figure 3 caption had a collapsed "_GoBack" bookmark sitting in
# the middle of the sentence, splitting it into two runs. Remove the
# bookmark and rejoin the sentence into a single run.
# ---------------------------------------------------------------------------
$rng3 = $d.Content.Duplicate
$found3 = $rng3.Find.Execute("images, matrix", $true, $false, $false, $false, $false, $true, 1, $false, "images, matrix", 2)
